$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of cell -> decimal literal (value prior to *1000 multiplication)
$cellFormulas = @{
    "B4"  = "0.004705"
    "E4"  = "0.0028219272"
    "H4"  = "0.0039976932"
    "B5"  = "0.0045644119"
    "E5"  = "0.0011199818"
    "H5"  = "0.0017048259"
    "B6"  = "0.0217759062"
    "E6"  = "0.0063845291"
    "H6"  = "0.0087333419"
    "B7"  = "0.0047509563"
    "E7"  = "0.0028627959"
    "H7"  = "0.00459951"
    "B8"  = "0.0035358382"
    "E8"  = "0.00239746"
    "H8"  = "0.0034524217"
    "B9"  = "0.0024074307"
    "E9"  = "0.0022949296"
    "H9"  = "0.0027690735"
    "B10" = "0.0018092166"
    "E10" = "0.0016388099"
    "H10" = "0.0018045769"
    "B15" = "0.006170377"
    "E15" = "0.002656911"
    "H15" = "0.0042586686"
    "B16" = "0.0022587786"
    "E16" = "0.0005842277"
    "H16" = "0.0042095557"
    "B17" = "0.0121152656"
    "E17" = "0.003853496"
    "H17" = "0.0230317658"
    "B18" = "0.0068265456"
    "E18" = "0.0030935058"
    "H18" = "0.0042487014"
    "B19" = "0.0054967424"
    "E19" = "0.0023973196"
    "H19" = "0.0030870698"
    "B20" = "0.0048682335"
    "E20" = "0.0021722161"
    "H20" = "0.002428469"
    "B21" = "0.0021708051"
    "E21" = "0.0020009328"
    "H21" = "0.0015388403"
}

foreach ($addr in $cellFormulas.Keys) {
    $lit = $cellFormulas[$addr]
    $ws.Range($addr).Formula = "=" + $lit + "*1000"
}

# Apply number format "0.000" to all cells that now hold the *1000 formulas
foreach ($addr in $cellFormulas.Keys) {
    $ws.Range($addr).NumberFormat = "0.000"
}

# Update the active selection cell to match the new selection
$ws.Range("H27").Select()
